$d = $word.ActiveDocument

$p23 = $d.Paragraphs.Item(23)
$findRng1 = $p23.Range
$findRng1.Find.Execute("con automóvil habilitado. No se puede", $true, $false, $false, $false, $false, $true, 1, $false, "con automóvil habilitado. ^pNo se puede", 2) | Out-Null
$p24 = $d.Paragraphs.Item(24)
$p24.Range.ListFormat.ListLevelNumber = 2

$findRng2 = $p24.Range
$findRng2.Find.Execute("un viaje existente. Un cliente", $true, $false, $false, $false, $false, $true, 1, $false, "un viaje existente. ^pUn cliente", 2) | Out-Null
$p25 = $d.Paragraphs.Item(25)
$p25.Range.ListFormat.ListLevelNumber = 2

$findRng3 = $p25.Range
$findRng3.Find.Execute("tiene un automóvil habilitado. La hora", $true, $false, $false, $false, $false, $true, 1, $false, "tiene un automóvil habilitado. ^pLa hora", 2) | Out-Null
$p26 = $d.Paragraphs.Item(26)
$p26.Range.ListFormat.ListLevelNumber = 2

$findRng4 = $p26.Range
$findRng4.Find.Execute("La fecha debe ser menor o igual a la fecha de hoy", $true, $false, $false, $false, $false, $true, 1, $false, "La fecha de inicio y la fecha de fin ingresadas deben ser menor o igual a la fecha y hora de hoy", 2) | Out-Null

$p26final = $d.Paragraphs.Item(26)
$dup = $p26final.Range.Duplicate
$dup.Collapse(0)  # wdCollapseEnd = 0
Write-Host "dup start/end after collapse end: " $dup.Start $dup.End
$dup.MoveEnd(1, -1) | Out-Null  # wdCharacter = 1, move end back 1 char
Write-Host "dup start/end after moveend -1: " $dup.Start $dup.End
$d.Bookmarks.Add("_GoBack", $dup) | Out-Null
$finalBm = $d.Bookmarks.Item("_GoBack")
Write-Host "finalBm: " $finalBm.Start $finalBm.End
